$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Change 1: "Led project to develop a hybrid streaming/batch pipeline, ..."
#   -> split into 3 runs, replacing ", " with " for social media advertising
#      data, " between "pipeline" and "shipping"
# ---------------------------------------------------------------------------

$r1 = $d.Content
$found = $r1.Find.Execute("pipeline, shipping", $true, $false, $false, $false, $false, $true, 1, $false, "pipeline for social media advertising data, shipping", 2)

$r1b = $d.Content
$target1 = "Led project to develop a hybrid streaming/batch pipeline for social media advertising data, shipping an MVP in only 7 weeks and scaling to tens of millions of "
$found1b = $r1b.Find.Execute($target1, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$start1 = $r1b.Start
$end1 = $r1b.End

$s1a = "Led project to develop a hybrid streaming/batch pipeline"
$s1b = " for social media advertising data, "

$seg1a = $d.Range($start1, $start1 + $s1a.Length)
$seg1a.Bold = 1
$seg1a.Bold = 0

$seg1b = $d.Range($start1 + $s1a.Length, $start1 + $s1a.Length + $s1b.Length)
$seg1b.Bold = 1
$seg1b.Bold = 0

$seg1c = $d.Range($start1 + $s1a.Length + $s1b.Length, $end1)
$seg1c.Bold = 1
$seg1c.Bold = 0

# ---------------------------------------------------------------------------
# Change 2: "Acted as the technical lead to deliver Amazon affiliate ..."
#   -> split into 3 runs, inserting "marketing " between "affiliate " and
#      "program"
# ---------------------------------------------------------------------------

$r2 = $d.Content
$found2 = $r2.Find.Execute("Amazon affiliate program", $true, $false, $false, $false, $false, $true, 1, $false, "Amazon affiliate marketing program", 2)

$r2b = $d.Content
$target2 = "Acted as the technical lead to deliver Amazon affiliate marketing program commerce analytics involving three teams and multiple contractors, providing experiment data to optimize conversions on a front page "
$found2b = $r2b.Find.Execute($target2, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$start2 = $r2b.Start
$end2 = $r2b.End

$s2a = "Acted as the technical lead to deliver Amazon affiliate "
$s2b = "marketing "

$seg2a = $d.Range($start2, $start2 + $s2a.Length)
$seg2a.Bold = 1
$seg2a.Bold = 0

$seg2b = $d.Range($start2 + $s2a.Length, $start2 + $s2a.Length + $s2b.Length)
$seg2b.Bold = 1
$seg2b.Bold = 0

$seg2c = $d.Range($start2 + $s2a.Length + $s2b.Length, $end2)
$seg2c.Bold = 1
$seg2c.Bold = 0

# ---------------------------------------------------------------------------
# Change 3: remove the "Document Languages: LaTeX, Markdown, Org,
#   ReStructuredText" paragraph entirely
# ---------------------------------------------------------------------------

$count = $d.Paragraphs.Count
for ($i = $count; $i -ge 1; $i--) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "Document Languages:*") {
        $p.Range.Delete()
        break
    }
}

Write-Host "Done"
